# "fix typos and other gremlins"
# Corrects three mis-typed / inconsistent PAF (95% CI) values in the IPF
# exposure table, and brings the sheet's cosmetic window/view state in
# line with where the author last left it (active cell + column widths).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Typo fixes: Pooled PAF (95% CI) column (D) ---------------------------
# "Wood dust" row: 6 (3-12) -> 4 (2-6)
$ws.Range("D4").Value = "4 (2-6)"
# "Agricultural dust" row: 12(0-40) -> 4 (0-12)  (also fixes missing space)
$ws.Range("D5").Value = "4 (0-12)"
# "Silica dust" row: 14 (8-20) -> 3 (2-5)
$ws.Range("D6").Value = "3 (2-5)"

# --- Window/view state: active cell moved from D13 to E9 ------------------
$ws.Range("E9").Select()

# --- Minor column width adjustments (as left by the editing session) ------
$ws.Columns.Item(1).ColumnWidth = 22
$ws.Columns.Item(2).ColumnWidth = 40.3333333333333
$ws.Columns.Item(3).ColumnWidth = 26.8333333333333
$ws.Columns.Item(4).ColumnWidth = 40.3333333333333
$ws.Columns.Item(5).ColumnWidth = 30.8333333333333

# --- Tab ratio of the split between sheet tabs and horizontal scrollbar ---
$excel.ActiveWindow.TabRatio = 0.993
